$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2").Value = "a.chagas@senff.com.br"
$ws.Range("C3").Select()
